# Apply updated odds values to row 6 of Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H6").Value = 3.5
$ws.Range("L6").Value = 5
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 9.5
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 1.85
$ws.Range("U6").Value = 1.83
$ws.Range("V6").Value = 1.83
$ws.Range("X6").Value = 8
$ws.Range("Z6").Value = 13
$ws.Range("AC6").Value = 9.5
$ws.Range("AD6").Value = 6.5
$ws.Range("AU6").Value = 8.5
$ws.Range("BB6").Value = 251
